$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new price text is numeric-looking so Excel keeps them as
# text (matching the source inlineStr cells) instead of auto-converting to numbers.
$ws.Range('D4:D6').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D10:D14').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D20:D32').NumberFormat = '@'
$ws.Range('D34:D42').NumberFormat = '@'
$ws.Range('D44:D46').NumberFormat = '@'
$ws.Range('D50:D51').NumberFormat = '@'

$ws.Range('D2').Value = '67.605.84'
$ws.Range('E2').Value = '  +3.65%  '
$ws.Range('D3').Value = '3.293.69'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '574.98'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('D6').Value = '177.18'
$ws.Range('E6').Value = '  -2.10%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D8').Value = '0.584'
$ws.Range('E8').Value = '  +3.51%  '
$ws.Range('D9').Value = '3.289.15'
$ws.Range('E9').Value = '  +0.53%  '
$ws.Range('D10').Value = '0.174'
$ws.Range('E10').Value = '  +0.62%  '
$ws.Range('D11').Value = '0.573'
$ws.Range('E11').Value = '  +1.80%  '
$ws.Range('D12').Value = '45.57'
$ws.Range('E12').Value = '  -0.40%  '
$ws.Range('D13').Value = '0.0000269'
$ws.Range('E13').Value = '  +3.20%  '
$ws.Range('D14').Value = '708.05'
$ws.Range('E14').Value = '  +14.98%  '
$ws.Range('D15').Value = '3.825.88'
$ws.Range('E15').Value = '  +1.09%  '
$ws.Range('D16').Value = '8.35'
$ws.Range('E16').Value = '  +0.51%  '
$ws.Range('D17').Value = '67.710.18'
$ws.Range('E17').Value = '  +3.59%  '
$ws.Range('E18').Value = '  +1.61%  '
$ws.Range('D19').Value = '3.304.35'
$ws.Range('E19').Value = '  +1.15%  '
$ws.Range('D20').Value = '17.41'
$ws.Range('E20').Value = '  -0.62%  '
$ws.Range('D21').Value = '10.76'
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('D22').Value = '0.891'
$ws.Range('E22').Value = '  +1.54%  '
$ws.Range('D23').Value = '16.81'
$ws.Range('E23').Value = '  -6.82%  '
$ws.Range('D24').Value = '5.15'
$ws.Range('E24').Value = '  +4.38%  '
$ws.Range('D25').Value = '98.87'
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('D26').Value = '3.92'
$ws.Range('E26').Value = '  +0.23%  '
$ws.Range('D27').Value = '2.72'
$ws.Range('E27').Value = '  +1.10%  '
$ws.Range('D28').Value = '9.31'
$ws.Range('E28').Value = '  +0.43%  '
$ws.Range('D29').Value = '33.02'
$ws.Range('E29').Value = '  +8.96%  '
$ws.Range('D30').Value = '8.43'
$ws.Range('E30').Value = '  +2.20%  '
$ws.Range('D31').Value = '6.66'
$ws.Range('E31').Value = '  +3.68%  '
$ws.Range('D32').Value = '578.22'
$ws.Range('E32').Value = '  +5.29%  '
$ws.Range('D33').Value = '3.900.78'
$ws.Range('E33').Value = '  +3.47%  '
$ws.Range('D34').Value = '10.82'
$ws.Range('E34').Value = '  +1.01%  '
$ws.Range('D35').Value = '0.103'
$ws.Range('E35').Value = '  +1.40%  '
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('D37').Value = '3.34'
$ws.Range('E37').Value = '  -5.72%  '
$ws.Range('D38').Value = '55.28'
$ws.Range('E38').Value = '  -0.44%  '
$ws.Range('D39').Value = '0.129'
$ws.Range('E39').Value = '  +2.96%  '
$ws.Range('D40').Value = '3.13'
$ws.Range('E40').Value = '  +1.14%  '
$ws.Range('D41').Value = '2.60'
$ws.Range('E41').Value = '  +2.21%  '
$ws.Range('D42').Value = '32.00'
$ws.Range('E42').Value = '  -0.37%  '
$ws.Range('D43').Value = '0.0₃0676'
$ws.Range('E43').Value = '  +1.48%  '
$ws.Range('D44').Value = '3.34'
$ws.Range('E44').Value = '  -1.81%  '
$ws.Range('D45').Value = '0.329'
$ws.Range('E45').Value = '  +1.03%  '
$ws.Range('D46').Value = '0.0409'
$ws.Range('E46').Value = '  +2.36%  '
$ws.Range('E47').Value = '  +2.48%  '
$ws.Range('E48').Value = '  +11.72%  '
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('D50').Value = '2.53'
$ws.Range('E50').Value = '  +2.31%  '
$ws.Range('D51').Value = '128.29'
$ws.Range('E51').Value = '  +0.29%  '
